$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.582.93"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.402.01"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.43"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.80"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +3.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.394.51"
$ws.Range("E8").Value = "  +2.15%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +15.11%  "
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.83"
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000282"
$ws.Range("E13").Value = "  +6.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.15"
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.947.54"
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.32"
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.404.71"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "65.527.49"
$ws.Range("E19").Value = "  +3.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.88"
$ws.Range("E20").Value = "  +2.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.993"
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "471.72"
$ws.Range("E22").Value = "  +16.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.15"
$ws.Range("E23").Value = "  +19.89%  "
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.68"
$ws.Range("E25").Value = "  +4.75%  "
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.90"
$ws.Range("E27").Value = "  +3.43%  "
$ws.Range("E28").Value = "  +6.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.89"
$ws.Range("E29").Value = "  +4.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.06"
$ws.Range("E30").Value = "  +7.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.71"
$ws.Range("E31").Value = "  +5.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.56"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.06"
$ws.Range("E33").Value = "  +10.36%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "582.03"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("E35").Value = "  +2.56%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.50"
$ws.Range("E38").Value = "  +3.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.75"
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("E40").Value = "  +2.51%  "
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.087.87"
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.22"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.48"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("E48").Value = "  +5.99%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.36"
$ws.Range("E50").Value = "  +4.88%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.28"
$ws.Range("E51").Value = "  +2.69%  "
